$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "11 Apr 2018 16:00:00.000'"
$ws.Range("C2").Value = 33
$ws.Range("D2").Value = -104

$ws.Columns("B:D").AutoFit() | Out-Null

$ws.Range("E2").Select() | Out-Null
